$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.959.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.521.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.517.41"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.81"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.579"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.34"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.076.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "616.55"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.935.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.508.24"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.886"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.76"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.61"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.83"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.61"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.91%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.51"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.45%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.35%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.93"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -8.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "623.71"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.75"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.13"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -13.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0445"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.383.83"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.326"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0$([char]0x2083)0695"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.10%  "
